$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell J1 = "Q8", matching style of existing header cells (I1)
$ws.Range("J1").Value = "Q8"
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2 values (B2:H2)
$ws.Range("B2").Value = -0.01758111206922311
$ws.Range("C2").Value = -0.9703086835434362
$ws.Range("D2").Value = -1.190204369659767
$ws.Range("E2").Value = 0.851657034662153
$ws.Range("F2").Value = 1.077585438625931
$ws.Range("G2").Value = -0.1996304584251192
$ws.Range("H2").Value = 0.1134005067055028

# Row 3 values (B3:G3)
$ws.Range("B3").Value = -0.8029421940374804
$ws.Range("C3").Value = -1.022837880153811
$ws.Range("D3").Value = 1.019023524168109
$ws.Range("E3").Value = 1.244951928131887
$ws.Range("F3").Value = -0.0322639689191633
$ws.Range("G3").Value = 0.2807669962114587

# Row 4 values (B4:J4) - new cells G4:J4 added
$ws.Range("B4").Value = -0.7936197797752114
$ws.Range("C4").Value = 1.248241624546709
$ws.Range("D4").Value = 1.474170028510487
$ws.Range("E4").Value = 0.1969541314594366
$ws.Range("F4").Value = 0.5099850965900585
$ws.Range("G4").Value = 0.007569982732279296
$ws.Range("H4").Value = 0.6592032301411037
$ws.Range("I4").Value = 0.4663036110991228
$ws.Range("J4").Value = -0.04443783748577212

# Row 5 values (B5:I5)
$ws.Range("B5").Value = 2.484163736993811
$ws.Range("C5").Value = 2.710092140957589
$ws.Range("D5").Value = 1.432876243906539
$ws.Range("E5").Value = 1.745907209037161
$ws.Range("F5").Value = 1.243492095179382
$ws.Range("G5").Value = 1.895125342588206
$ws.Range("H5").Value = 1.702225723546225
$ws.Range("I5").Value = 1.19148427496133

# Row 6 values (B6:H6)
$ws.Range("B6").Value = 1.371380565536508
$ws.Range("C6").Value = 0.09416466848545757
$ws.Range("D6").Value = 0.4071956336160796
$ws.Range("E6").Value = -0.09521948024169971
$ws.Range("F6").Value = 0.5564137671671248
$ws.Range("G6").Value = 0.3635141481251438
$ws.Range("H6").Value = -0.1472273004597511

# Row 7 values (B7:G7)
$ws.Range("B7").Value = 0.2659007569564139
$ws.Range("C7").Value = 0.5789317220870359
$ws.Range("D7").Value = 0.07651660822925663
$ws.Range("E7").Value = 0.7281498556380811
$ws.Range("F7").Value = 0.5352502365961002
$ws.Range("G7").Value = 0.02450878801120521

# Row 8 values (B8:I8) - new cells G8:I8 added
$ws.Range("B8").Value = 1.260690851164143
$ws.Range("C8").Value = 0.7582757373063643
$ws.Range("D8").Value = 1.409908984715189
$ws.Range("E8").Value = 1.217009365673208
$ws.Range("F8").Value = 0.7062679170883128
$ws.Range("G8").Value = 0.911668649685511
$ws.Range("H8").Value = 0.5051827077222001
$ws.Range("I8").Value = 0.9402868649905415

# Row 9 values (B9:H9)
$ws.Range("B9").Value = 0.2946970959196917
$ws.Range("C9").Value = 0.9463303433285162
$ws.Range("D9").Value = 0.7534307242865352
$ws.Range("E9").Value = 0.2426892757016403
$ws.Range("F9").Value = 0.4480900082988384
$ws.Range("G9").Value = 0.04160406633552749
$ws.Range("H9").Value = 0.4767082236038689

# Row 10 values (B10:G10)
$ws.Range("B10").Value = 0.3856725119803543
$ws.Range("C10").Value = 0.1927728929383733
$ws.Range("D10").Value = -0.3179685556465216
$ws.Range("E10").Value = -0.1125678230493235
$ws.Range("F10").Value = -0.5190537650126344
$ws.Range("G10").Value = -0.08394960774429301

# Row 11 values (B11:F11)
$ws.Range("B11").Value = 0.2952882579329085
$ws.Range("C11").Value = -0.2154531906519864
$ws.Range("D11").Value = -0.01005245805478834
$ws.Range("E11").Value = -0.4165384000180992
$ws.Range("F11").Value = 0.01856575725024216

# Row 12 values (B12:E12)
$ws.Range("B12").Value = -0.2970557949068323
$ws.Range("C12").Value = -0.09165506230963413
$ws.Range("D12").Value = -0.4981410042729451
$ws.Range("E12").Value = -0.06303684700460363

# Row 13 values (B13:D13)
$ws.Range("B13").Value = -0.1069508448768545
$ws.Range("C13").Value = -0.5134367868401654
$ws.Range("D13").Value = -0.07833262957182399

# Row 14 values (B14:C14)
$ws.Range("B14").Value = -0.6053253388254292
$ws.Range("C14").Value = -0.1702211815570877

# Row 15 values (B15)
$ws.Range("B15").Value = 0.06843616378760228

# Row 16: no numeric value cells, only the label in column A (unchanged)
